# Auto-generated edit script
# Applies updated currentAveragePrice / LevePrice / LeveProfit values
# to the Chocobo Profits workbook, per sheet (job class).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 29990
$ws.Range("J10").Value = 29990
$ws.Range("L10").Value = 29990
$ws.Range("N10").Value = -30576
$ws.Range("H28").Value = 620.05554
$ws.Range("I28").Value = 512.38464
$ws.Range("J28").Value = 900
$ws.Range("K28").Value = 512.38464
$ws.Range("L28").Value = 900
$ws.Range("M28").Value = -27.38463999999999
$ws.Range("N28").Value = -1870
$ws.Range("H64").Value = 3773.4614
$ws.Range("I64").Value = 3500
$ws.Range("J64").Value = 5277.5
$ws.Range("K64").Value = 3500
$ws.Range("L64").Value = 5277.5
$ws.Range("M64").Value = -3252
$ws.Range("N64").Value = -5773.5
$ws.Range("H67").Value = 3773.4614
$ws.Range("I67").Value = 3500
$ws.Range("J67").Value = 5277.5
$ws.Range("K67").Value = 3500
$ws.Range("L67").Value = 5277.5
$ws.Range("M67").Value = -2642
$ws.Range("N67").Value = -6993.5
$ws.Range("H74").Value = 50004900
$ws.Range("J74").Value = 9800
$ws.Range("L74").Value = 9800
$ws.Range("N74").Value = -11672
$ws.Range("H76").Value = 3187.375
$ws.Range("I76").Value = 3195.3635
$ws.Range("J76").Value = 3099.5
$ws.Range("K76").Value = 3195.3635
$ws.Range("L76").Value = 3099.5
$ws.Range("M76").Value = -2880.3635
$ws.Range("N76").Value = -3729.5
$ws.Range("H77").Value = 50004900
$ws.Range("J77").Value = 9800
$ws.Range("L77").Value = 49000
$ws.Range("N77").Value = -58360
$ws.Range("H79").Value = 3187.375
$ws.Range("I79").Value = 3195.3635
$ws.Range("J79").Value = 3099.5
$ws.Range("K79").Value = 3195.3635
$ws.Range("L79").Value = 3099.5
$ws.Range("M79").Value = -2103.3635
$ws.Range("N79").Value = -5283.5
$ws.Range("H107").Value = 1882.2142
$ws.Range("I107").Value = 2321.6667
$ws.Range("J107").Value = 1091.2
$ws.Range("K107").Value = 2321.6667
$ws.Range("L107").Value = 1091.2
$ws.Range("M107").Value = -401.6667000000002
$ws.Range("N107").Value = -4931.2
$ws.Range("H131").Value = 2667.3572
$ws.Range("I131").Value = 1355.375
$ws.Range("J131").Value = 4416.6665
$ws.Range("K131").Value = 4066.125
$ws.Range("L131").Value = 13249.9995
$ws.Range("M131").Value = 973.875
$ws.Range("N131").Value = -23329.9995
$ws.Range("H133").Value = 51220
$ws.Range("J133").Value = 51220
$ws.Range("L133").Value = 51220
$ws.Range("N133").Value = -61340

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3349.8918
$ws.Range("I32").Value = 2983.7231
$ws.Range("J32").Value = 5994.4443
$ws.Range("K32").Value = 2983.7231
$ws.Range("L32").Value = 5994.4443
$ws.Range("M32").Value = -2696.7231
$ws.Range("N32").Value = -6568.4443
$ws.Range("H88").Value = 9527481
$ws.Range("I88").Value = 13336833
$ws.Range("J88").Value = 4100
$ws.Range("K88").Value = 13336833
$ws.Range("L88").Value = 4100
$ws.Range("M88").Value = -13336427
$ws.Range("N88").Value = -4912
$ws.Range("H91").Value = 9527481
$ws.Range("I91").Value = 13336833
$ws.Range("J91").Value = 4100
$ws.Range("K91").Value = 13336833
$ws.Range("L91").Value = 4100
$ws.Range("M91").Value = -13335429
$ws.Range("N91").Value = -6908
$ws.Range("H132").Value = 2793.652
$ws.Range("I132").Value = 1417.5385
$ws.Range("J132").Value = 4582.6
$ws.Range("K132").Value = 4252.6155
$ws.Range("L132").Value = 13747.8
$ws.Range("M132").Value = -1722.6155
$ws.Range("N132").Value = -18807.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2250.0625
$ws.Range("I86").Value = 2059.5715
$ws.Range("J86").Value = 2398.2222
$ws.Range("K86").Value = 2059.5715
$ws.Range("L86").Value = 2398.2222
$ws.Range("M86").Value = -936.5715
$ws.Range("N86").Value = -4644.2222
$ws.Range("H89").Value = 2250.0625
$ws.Range("I89").Value = 2059.5715
$ws.Range("J89").Value = 2398.2222
$ws.Range("K89").Value = 10297.8575
$ws.Range("L89").Value = 11991.111
$ws.Range("M89").Value = -4681.8575
$ws.Range("N89").Value = -23223.111
$ws.Range("H95").Value = 36700
$ws.Range("J95").Value = 36700
$ws.Range("L95").Value = 36700
$ws.Range("N95").Value = -42192
$ws.Range("H134").Value = 2576.2
$ws.Range("I134").Value = 1696.2222
$ws.Range("K134").Value = 5088.6666
$ws.Range("M134").Value = -2553.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2451.742
$ws.Range("I31").Value = 1003.5625
$ws.Range("J31").Value = 3996.4666
$ws.Range("K31").Value = 1003.5625
$ws.Range("L31").Value = 3996.4666
$ws.Range("M31").Value = -708.5625
$ws.Range("N31").Value = -4586.4666
$ws.Range("H34").Value = 2451.742
$ws.Range("I34").Value = 1003.5625
$ws.Range("J34").Value = 3996.4666
$ws.Range("K34").Value = 1003.5625
$ws.Range("L34").Value = 3996.4666
$ws.Range("M34").Value = -801.5625
$ws.Range("N34").Value = -4400.4666
$ws.Range("H62").Value = 5150
$ws.Range("J62").Value = 7500
$ws.Range("L62").Value = 7500
$ws.Range("N62").Value = -8748
$ws.Range("H65").Value = 5150
$ws.Range("J65").Value = 7500
$ws.Range("L65").Value = 37500
$ws.Range("N65").Value = -43740
$ws.Range("H134").Value = 3254.7778
$ws.Range("I134").Value = 950
$ws.Range("J134").Value = 3913.2856
$ws.Range("K134").Value = 2850
$ws.Range("L134").Value = 11739.8568
$ws.Range("M134").Value = -315
$ws.Range("N134").Value = -16809.8568
$ws.Range("H137").Value = 44593.332
$ws.Range("J137").Value = 44593.332
$ws.Range("L137").Value = 44593.332
$ws.Range("N137").Value = -54793.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 1552.3846
$ws.Range("I25").Value = 2250.5
$ws.Range("K25").Value = 6751.5
$ws.Range("M25").Value = -6582.5
$ws.Range("H30").Value = 1552.3846
$ws.Range("I30").Value = 2250.5
$ws.Range("K30").Value = 6751.5
$ws.Range("M30").Value = -6649.5
$ws.Range("H113").Value = 713.0732
$ws.Range("I113").Value = 708.37036
$ws.Range("J113").Value = 722.1429000000001
$ws.Range("K113").Value = 2125.11108
$ws.Range("L113").Value = 2166.4287
$ws.Range("M113").Value = 44.88891999999987
$ws.Range("N113").Value = -6506.4287
$ws.Range("H131").Value = 13158762
$ws.Range("J131").Value = 905.5833
$ws.Range("L131").Value = 2716.7499
$ws.Range("N131").Value = -12796.7499

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6176.6816
$ws.Range("I70").Value = 5767.1294
$ws.Range("J70").Value = 8019.6665
$ws.Range("K70").Value = 5767.1294
$ws.Range("L70").Value = 8019.6665
$ws.Range("M70").Value = -5497.1294
$ws.Range("N70").Value = -8559.666499999999
$ws.Range("H73").Value = 6176.6816
$ws.Range("I73").Value = 5767.1294
$ws.Range("J73").Value = 8019.6665
$ws.Range("K73").Value = 5767.1294
$ws.Range("L73").Value = 8019.6665
$ws.Range("M73").Value = -4831.1294
$ws.Range("N73").Value = -9891.666499999999
$ws.Range("H113").Value = 1700
$ws.Range("J113").Value = 1700
$ws.Range("L113").Value = 1700
$ws.Range("N113").Value = -6040
$ws.Range("H126").Value = 3329.45
$ws.Range("I126").Value = 2976.218
$ws.Range("J126").Value = 4581.8184
$ws.Range("K126").Value = 8928.653999999999
$ws.Range("L126").Value = 13745.4552
$ws.Range("M126").Value = -6458.653999999999
$ws.Range("N126").Value = -18685.4552
$ws.Range("H132").Value = 2495.0938
$ws.Range("I132").Value = 2042.3572
$ws.Range("J132").Value = 2847.2222
$ws.Range("K132").Value = 6127.071599999999
$ws.Range("L132").Value = 8541.6666
$ws.Range("M132").Value = -3597.071599999999
$ws.Range("N132").Value = -13601.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3708.0588
$ws.Range("I7").Value = 2056
$ws.Range("J7").Value = 6068.143
$ws.Range("K7").Value = 2056
$ws.Range("L7").Value = 6068.143
$ws.Range("M7").Value = -1944
$ws.Range("N7").Value = -6292.143
$ws.Range("H35").Value = 300
$ws.Range("I35").Value = 300
$ws.Range("K35").Value = 300
$ws.Range("M35").Value = 36
$ws.Range("H126").Value = 3708.0588
$ws.Range("I126").Value = 2056
$ws.Range("J126").Value = 6068.143
$ws.Range("K126").Value = 6168
$ws.Range("L126").Value = 18204.429
$ws.Range("M126").Value = -3698
$ws.Range("N126").Value = -23144.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 27516.666
$ws.Range("J92").Value = 27516.666
$ws.Range("L92").Value = 27516.666
$ws.Range("N92").Value = -32508.666
$ws.Range("H122").Value = 5446.385
$ws.Range("I122").Value = 3556.8572
$ws.Range("J122").Value = 7650.8335
$ws.Range("K122").Value = 10670.5716
$ws.Range("L122").Value = 22952.5005
$ws.Range("M122").Value = -8220.571599999999
$ws.Range("N122").Value = -27852.5005
$ws.Range("H136").Value = 4303.9
$ws.Range("I136").Value = 1873.3334
$ws.Range("K136").Value = 5620.0002
$ws.Range("M136").Value = -3070.0002

Write-Output "Applied 232 cell updates across 8 sheets."